$d = $word.ActiveDocument

# Locate the "1PE" Heading 2 paragraph, then remove the very next paragraph,
# which is the italic subtitle paragraph ("Первое послание Петра") that
# duplicates the book title shown later in the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "1PE" -and $p.Style.NameLocal -eq "Heading 2") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $nextRange = $target.Next().Range
    $nextRange.Delete()
}
